# Updated capital structure database
# Apply updated metrics for Papua New Guinea - Metals & Mining rows (2 and 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.04849999999999999
$ws.Range("G2").Value = -1.972727272727272
$ws.Range("H2").Value = -1.972727272727272
$ws.Range("I2").Value = -2.218181818181818
$ws.Range("J2").Value = -2.218181818181818
$ws.Range("K2").Value = -2.64
$ws.Range("L2").Value = -2.4
$ws.Range("U2").Value = 1.41
$ws.Range("V2").Value = 0.0157190635451505
$ws.Range("W2").Value = -0.08627450980392157
$ws.Range("X2").Value = 0.09283888990821114
$ws.Range("Y2").Value = -0.1791133997121327
$ws.Range("Z2").Value = 0.03566218187712758
$ws.Range("AA2").Value = -0.07910520343653753
$ws.Range("AB2").Value = 0.09269318667149577
$ws.Range("AC2").Value = -0.1717983901080333
$ws.Range("AD2").Value = 0.347
$ws.Range("AF2").Value = 0.347
$ws.Range("AG2").Value = -1.063
$ws.Range("AH2").Value = 0.003853543149688496
$ws.Range("AI2").Value = 0.01322055854002362
$ws.Range("AJ2").Value = -0.01199273441113756
$ws.Range("AK2").Value = -0.04279904980472682
$ws.Range("AL2").Value = 0.044
$ws.Range("AM2").Value = 0.044
$ws.Range("AN2").Value = -0.151528384279476
$ws.Range("AO2").Value = -55.45454545454545
$ws.Range("AP2").Value = 0.4641921397379912
$ws.Range("AQ2").Value = -55.45454545454545

# Row 3
$ws.Range("D3").Value = -0.04849999999999999
$ws.Range("G3").Value = -1.972727272727272
$ws.Range("H3").Value = -1.972727272727272
$ws.Range("I3").Value = -2.218181818181818
$ws.Range("J3").Value = -2.218181818181818
$ws.Range("K3").Value = -2.64
$ws.Range("L3").Value = -2.4
$ws.Range("U3").Value = 1.41
$ws.Range("V3").Value = 0.0157190635451505
$ws.Range("W3").Value = -0.08627450980392157
$ws.Range("X3").Value = 0.09283888990821114
$ws.Range("Y3").Value = -0.1791133997121327
$ws.Range("Z3").Value = 0.03566218187712758
$ws.Range("AA3").Value = -0.07910520343653753
$ws.Range("AB3").Value = 0.09269318667149577
$ws.Range("AC3").Value = -0.1717983901080333
$ws.Range("AD3").Value = 0.347
$ws.Range("AF3").Value = 0.347
$ws.Range("AG3").Value = -1.063
$ws.Range("AH3").Value = 0.003853543149688496
$ws.Range("AI3").Value = 0.01322055854002362
$ws.Range("AJ3").Value = -0.01199273441113756
$ws.Range("AK3").Value = -0.04279904980472682
$ws.Range("AL3").Value = 0.044
$ws.Range("AM3").Value = 0.044
$ws.Range("AN3").Value = -0.151528384279476
$ws.Range("AO3").Value = -55.45454545454545
$ws.Range("AP3").Value = 0.4641921397379912
$ws.Range("AQ3").Value = -55.45454545454545
